# Insert a new data row before row 128 (the existing rows 128-207 shift down to 129-208).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(128).Insert()

$ws.Cells.Item(128, 1).Value = 4
$ws.Cells.Item(128, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(128, 3).Value = "Los Lagos"
$ws.Cells.Item(128, 4).Value = "2022-02-04"
$ws.Cells.Item(128, 5).Value = 10
$ws.Cells.Item(128, 6).Value = 100112044
$ws.Cells.Item(128, 7).Value = "Perejil"
$ws.Cells.Item(128, 8).Value = "Sin especificar"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 140
$ws.Cells.Item(128, 11).Value = 5000
$ws.Cells.Item(128, 12).Value = 5000
$ws.Cells.Item(128, 13).Value = 5000
$ws.Cells.Item(128, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(128, 15).Value = "Región Metropolitana"
$ws.Cells.Item(128, 16).Value = 1667
$ws.Cells.Item(128, 17).Value = 3
$ws.Cells.Item(128, 18).Value = "Hortaliza"
